# Telecom-Media.xlsx update — append newly tracked tickers to the "Main"
# sheet and refresh the view zoom, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New name/ticker rows to append after the existing data (rows 3-33).
$newRows = @(
    @("Singapore Telecom", "ST SP"),
    @("Telstra", "TLS AU"),
    @("Telus", "T CN"),
    @("Cellnex", "CLNX"),
    @("Orange", "ORA FP"),
    @("Swisscom", "SCMN SW"),
    @("Telekom Indonesia", "TLKM IJ"),
    @("Telefonica", "TEF SM"),
    @("Wolters Kluwer", "WKL NA"),
    @("Sirius XM", "SIRI")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $name = $newRows[$i][0]
    $ticker = $newRows[$i][1]

    $ws.Range("A$row").Value = "x"
    $ws.Range("B$row").Value = $name
    $ws.Range("C$row").Value = $ticker
}

# Match the author's updated zoom level on the sheet view.
$excel.ActiveWindow.Zoom = 175
